$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append new rows 245:254 with quarterly data (dates in A, hours in B) ---
$dates  = @(43525, 43556, 43586, 43617, 43647, 43678, 43709, 43739, 43770, 43800)
$values = @(36.5,   37.4,  38,    37.7,  37.4,  37.8,  36.2,  36,    35.5,  36.8)

for ($i = 0; $i -lt $dates.Length; $i++) {
    $row = 245 + $i
    $ws.Cells.Item($row, 1).Value = $dates[$i]
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Style column A245:A254 like A244 (date format) by copying format from the existing cell
$ws.Range("A244").Copy()
$ws.Range("A245:A254").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# Build the new style (number + right/center alignment + thin left/right border) on B245 first
$b1 = $ws.Range("B245")
$b1.NumberFormat = "#,##0.00"
$b1.HorizontalAlignment = -4152   # xlRight
$b1.VerticalAlignment = -4108     # xlCenter
$b1.Borders.Item(7).LineStyle = 1   # xlEdgeLeft
$b1.Borders.Item(10).LineStyle = 1  # xlEdgeRight

# Propagate the same style to the rest of the column via format copy (keeps a single shared style)
$b1.Copy()
$ws.Range("B246:B254").PasteSpecial(-4122) # xlPasteFormats
$excel.CutCopyMode = $false

# --- Update D94:D97 to use quarterly-average formulas referencing the new rows ---
$ws.Range("D94").Formula = "=(B243+B244+B245)/3"
$ws.Range("D95").Formula = "=(B246+B247+B248)/3"
$ws.Range("D96").Formula = "=(B249+B250+B251)/3"
$ws.Range("D97").Formula = "=(B252+B253+B254)/3"

# --- Update the active selection to D98 (matches author's final cursor position) ---
$ws.Range("D98").Select()
